# Intake sheet: turn the blank "CASE NO._____" + "___" signature line into
# an underlined "CASE NO. _${case_no}_" merge field, and drop a _GoBack
# bookmark right after the trio of tab stops that follows it.

$d = $word.ActiveDocument

# Locate "CASE NO." followed by the 8 underscores (two runs: "_____" + "___").
# NOTE: keep a single Range handle and mutate it via Find -- re-fetching
# $d.Content after the call yields a fresh, unmoved full-document Range.
$hit = $d.Content
$found = $hit.Find.Execute("CASE NO.________", $true, $false, $false,
                            $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'CASE NO.________' placeholder text."
}
$start = $hit.Start

# Replace the 16-character run with the new 21-character field text.
$full = $d.Range($start, $start + 16)
$full.Text = "CASE NO. _`${case_no}_"

# Force a run split between "CASE NO." and the following space (both stay
# bold) by toggling Bold off/on over just the space -- Word/iron_native
# merges runs with identical formatting, so a no-op format write on a
# sub-range is what keeps the boundary.
$spaceRange = $d.Range($start + 8, $start + 9)
$spaceRange.Font.Bold = 0
$spaceRange.Font.Bold = 1

# Underline "_${case_no}_" (the leading underscore, the field, and the
# trailing underscore).
$underlineRange = $d.Range($start + 9, $start + 21)
$underlineRange.Font.Underline = 1

# The merge field placeholder itself is underlined but not bold.
$placeholderRange = $d.Range($start + 10, $start + 20)
$placeholderRange.Font.Bold = 0

# The three existing <w:tab/> runs survive untouched right after our new
# text; drop a collapsed _GoBack bookmark immediately after them (Word
# stamps this on every manual edit/save).
$bookmarkPos = $start + 21 + 3
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
